$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new value in E1 with accounting/currency number format (numFmtId 44)
$ws.Range("E1").Value = 100.75
$ws.Range("E1").NumberFormat = '_("$"* #,##0.00_);_("$"* \(#,##0.00\);_("$"* "-"??_);_(@_)'

# Update selection to E1 (matches final selection in the diff)
$ws.Range("E1").Select()
